# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; B = 'Bitcoin'; C = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D = '30.320.86'; E = '  -0.07%  '; ForceText = $false; UpdateBC = $false },
    @{ Row = 3; B = 'Ethereum'; C = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D = '1.884.14'; E = '  -0.95%  '; ForceText = $false; UpdateBC = $false },
    @{ Row = 4; B = 'TetherUSD'; C = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D = '0.9981'; E = '  -0.43%  '; ForceText = $true; UpdateBC = $false },
    @{ Row = 5; B = 'BNB'; C = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D = '237.89'; E = '  +0.18%  '; ForceText = $true; UpdateBC = $false },
    @{ Row = 6; B = 'USDC'; C = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D = '0.9988'; E = '  -0.32%  '; ForceText = $true; UpdateBC = $false },
    @{ Row = 7; B = 'XRP'; C = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D = '0.4673'; E = '  -0.84%  '; ForceText = $true; UpdateBC = $false },
    @{ Row = 8; B = 'Cardano'; C = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D = '0.2813'; E = '  -0.02%  '; ForceText = $true; UpdateBC = $false },
    @{ Row = 9; B = 'Dogecoin'; C = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D = '0.06554'; E = '  -0.59%  '; ForceText = $true; UpdateBC = $false },
    @{ Row = 10; B = 'Solana'; C = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D = '19.58'; E = '  +5.60%  '; ForceText = $true; UpdateBC = $false },
    @{ Row = 11; B = 'Litecoin'; C = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D = '98.50'; E = '  -0.84%  '; ForceText = $true; UpdateBC = $false },
    @{ Row = 12; B = 'TRON'; C = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D = '0.07731'; E = '  +0.29%  '; ForceText = $true; UpdateBC = $false },
    @{ Row = 13; B = 'WrappedEther'; C = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D = '1.883.32'; E = '  -0.99%  '; ForceText = $false; UpdateBC = $false },
    @{ Row = 14; B = 'Polkadot'; C = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D = '5.125'; E = '  -0.22%  '; ForceText = $true; UpdateBC = $false },
    @{ Row = 15; B = 'Polygon'; C = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D = '0.6686'; E = '  +0.97%  '; ForceText = $true; UpdateBC = $false },
    @{ Row = 16; B = 'BitcoinCash'; C = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D = '285.64'; E = '  +13.67%  '; ForceText = $true; UpdateBC = $false },
    @{ Row = 17; B = 'WrappedBTC'; C = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D = '30.305.07'; E = '  -0.22%  '; ForceText = $false; UpdateBC = $false },
    @{ Row = 18; B = 'Dai'; C = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; D = '0.9992'; E = '  -0.30%  '; ForceText = $true; UpdateBC = $false },
    @{ Row = 19; B = 'WrappedliquidstakedEther2.0'; C = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D = '2.132.05'; E = '  -1.15%  '; ForceText = $false; UpdateBC = $false },
    @{ Row = 20; B = 'Avalanche'; C = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D = '12.56'; E = '  +0.13%  '; ForceText = $true; UpdateBC = $false },
    @{ Row = 21; B = 'ShibaInu'; C = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D = '0.000007291'; E = '  -1.34%  '; ForceText = $true; UpdateBC = $false },
    @{ Row = 22; B = 'Uniswap'; C = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D = '5.310'; E = '  -0.58%  '; ForceText = $true; UpdateBC = $false },
    @{ Row = 23; B = 'BinanceUSD'; C = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D = '0.9981'; E = '  -0.49%  '; ForceText = $true; UpdateBC = $false },
    @{ Row = 24; B = 'BitDAO'; C = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'; D = '0.4574'; E = '  +0.17%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 25; B = 'Chainlink'; C = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D = '6.191'; E = '  -0.89%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 26; B = 'Monero'; C = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D = '167.47'; E = '  +2.10%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 27; B = 'Cosmos'; C = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D = '9.264'; E = '  -0.27%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 28; B = 'EthereumClassic'; C = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D = '19.05'; E = '  +1.87%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 29; B = 'LidoDAOToken'; C = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D = '1.985'; E = '  -2.29%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 30; B = 'Toncoin'; C = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D = '1.370'; E = '  -0.46%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 31; B = 'Stellar'; C = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D = '0.09848'; E = '  -2.06%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 32; B = 'Filecoin'; C = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D = '4.461'; E = '  -3.00%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 33; B = 'PancakeSwap'; C = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D = '1.492'; E = '  -0.86%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 34; B = 'InternetComputer(DFINITY)'; C = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D = '4.190'; E = '  -0.27%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 35; B = 'Hedera'; C = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; D = '0.04673'; E = '  -0.50%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 36; B = 'ImmutableX'; C = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D = '0.7089'; E = '  -1.56%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 37; B = 'ARBITRUM'; C = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D = '1.096'; E = '  -0.27%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 38; B = 'Frax'; C = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'; D = '0.9980'; E = '  -0.33%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 39; B = 'HuobiToken'; C = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D = '2.703'; E = '  -0.55%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 40; B = 'VeChain'; C = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D = '0.01871'; E = '  -1.32%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 41; B = 'FraxShare'; C = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D = '6.676'; E = '  +7.70%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 42; B = 'MXToken'; C = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D = '2.517'; E = '  -3.01%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 43; B = 'Aave'; C = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D = '72.30'; E = '  +0.12%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 44; B = 'TrustWalletToken'; C = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D = '0.8688'; E = '  +1.67%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 45; B = 'RenderToken'; C = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D = '1.961'; E = '  -0.12%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 46; B = 'Quant'; C = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D = '103.96'; E = '  -1.99%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 47; B = 'PaxDollar'; C = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D = '0.9983'; E = '  -0.33%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 48; B = 'TheSandbox'; C = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D = '0.4188'; E = '  -0.05%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 49; B = 'Maker'; C = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D = '992.70'; E = '  -4.09%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 50; B = 'Aptos'; C = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D = '7.246'; E = '  -1.60%  '; ForceText = $true; UpdateBC = $true },
    @{ Row = 51; B = 'EnergySwap'; C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D = '9.376'; E = '  +7.80%  '; ForceText = $true; UpdateBC = $true }
)

foreach ($item in $rows) {
    if ($item.UpdateBC) {
        $ws.Cells.Item($item.Row, 2).Value = $item.B
        $ws.Cells.Item($item.Row, 3).Value = $item.C
    }

    $dCell = $ws.Cells.Item($item.Row, 4)
    if ($item.ForceText) {
        $dCell.NumberFormat = "@"
    }
    $dCell.Value = $item.D

    $eCell = $ws.Cells.Item($item.Row, 5)
    $eCell.Value = $item.E
}
